# edit.ps1 - Applies the OOXML diff to Installation.docx via Word COM interop
# Summary of changes:
#  1. Merge the three split runs "Ubuntu 2" / "0" / ".04 vorgesehen." in the
#     "Installation von Garden" section into a single run (text unchanged).
#  2. Insert <w:lastRenderedPageBreak/> into the run holding the final
#     "Textfeld 4" code-block drawing (end of document body).
#  3. Remove the "INTERNAL" classification text-box (the whole <w:r> holding
#     the drawing/VML fallback) from every footer (footer1/2/3.xml), leaving
#     an empty "Fuzeile"-styled paragraph in each footer.

$d = $word.ActiveDocument

# --- 1) Fix the duplicated "Garden" installation paragraph: merge runs ----
# The original text is split across three runs ("Ubuntu 2" + "0" + ".04 ...").
# A plain Find&Replace over the full sentence lets Word re-emit it as a
# single run while leaving the visible text identical.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$oldText = "Die Installation von Harmonic ist für Ubuntu 20.04 vorgesehen."
$newText = "Die Installation von Harmonic ist für Ubuntu 20.04 vorgesehen."
[void]$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# --- 2) Insert <w:lastRenderedPageBreak/> before the last drawing's --------
# mc:AlternateContent, inside its existing run (right after <w:rPr>).
# We locate the paragraph by its embedded InlineShape (the last one in the
# document - "Textfeld 4") and replace its content with the same content
# plus the lastRenderedPageBreak marker, via InsertXML on that Range only.
$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.InlineShapes.Count -gt 0) {
        $targetParagraph = $p
    }
}
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w14:paraId="3C7D1B43" w14:textId="77777777" w:rsidR="00032632" w:rsidRDefault="00032632" w:rsidP="00032632"><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="4AA2BED8" wp14:editId="6080311E"><wp:extent cx="5715000" cy="1404620"/><wp:effectExtent l="0" t="0" r="19050" b="13335"/><wp:docPr id="4" name="Textfeld 4"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"><a:spLocks noChangeArrowheads="1"/></wps:cNvSpPr><wps:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5715000" cy="1404620"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFFFFF"/></a:solidFill><a:ln w="9525"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></wps:spPr><wps:txbx><w:txbxContent><w:p w14:paraId="6BFD21DC" w14:textId="77777777" w:rsidR="00032632" w:rsidRDefault="00032632" w:rsidP="00032632"><w:pPr><w:pStyle w:val="HTMLVorformatiert"/><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>apt-get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> update</w:t></w:r></w:p><w:p w14:paraId="6908EB02" w14:textId="77777777" w:rsidR="00032632" w:rsidRPr="00032632" w:rsidRDefault="00032632" w:rsidP="00032632"><w:pPr><w:pStyle w:val="HTMLVorformatiert"/><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr></w:p><w:p w14:paraId="162C8FC2" w14:textId="0CF419CE" w:rsidR="00032632" w:rsidRPr="00032632" w:rsidRDefault="00032632" w:rsidP="00032632"><w:pPr><w:pStyle w:val="HTMLVorformatiert"/><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>apt-get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>install</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>gz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>garden</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" anchor="t" anchorCtr="0"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="4AA2BED8" id="Textfeld 4" o:spid="_x0000_s1030" type="#_x0000_t202" style="width:450pt;height:110.6pt;visibility:visible;mso-wrap-style:square;mso-left-percent:-10001;mso-top-percent:-10001;mso-position-horizontal:absolute;mso-position-horizontal-relative:char;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-left-percent:-10001;mso-top-percent:-10001;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQBLbB3AGAIAACcEAAAOAAAAZHJzL2Uyb0RvYy54bWysU9tu2zAMfR+wfxD0vtgOkl6MOkWXLsOA&#10;7gJ0+wBZlmNhsqhRSuzs60fJaZp12MswPwiiSR0eHpI3t2Nv2F6h12ArXsxyzpSV0Gi7rfi3r5s3&#10;V5z5IGwjDFhV8YPy/Hb1+tXN4Eo1hw5Mo5ARiPXl4CreheDKLPOyU73wM3DKkrMF7EUgE7dZg2Ig&#10;9N5k8zy/yAbAxiFI5T39vZ+cfJXw21bJ8LltvQrMVJy4hXRiOut4ZqsbUW5RuE7LIw3xDyx6oS0l&#10;PUHdiyDYDvUfUL2WCB7aMJPQZ9C2WqpUA1VT5C+qeeyEU6kWEse7k0z+/8HKT/tH9wVZGN/CSA1M&#10;RXj3APK7ZxbWnbBbdYcIQ6dEQ4mLKFk2OF8en0apfekjSD18hIaaLHYBEtDYYh9VoToZoVMDDifR&#10;1RiYpJ/Ly2KZ5+SS5CsW+eJintqSifLpuUMf3ivoWbxUHKmrCV7sH3yIdET5FBKzeTC62WhjkoHb&#10;em2Q7QVNwCZ9qYIXYcayoeLXy/lyUuCvEEQ1sp2y/pap14FG2ei+4lenIFFG3d7ZJg1aENpMd6Js&#10;7FHIqN2kYhjrkemm4ouYIOpaQ3MgZRGmyaVNo0sH+JOzgaa24v7HTqDizHyw1J3rYrGIY56MxfKS&#10;pGR47qnPPcJKgqp44Gy6rkNajaSbu6MubnTS95nJkTJNY5L9uDlx3M/tFPW836tfAAAA//8DAFBL&#10;AwQUAAYACAAAACEAlBtql9sAAAAFAQAADwAAAGRycy9kb3ducmV2LnhtbEyPwU7DMBBE70j8g7WV&#10;uFGnkUAQ4lQVVc+UUglxc+xtHDVeh9hNU76ehQtcRhrNauZtuZx8J0YcYhtIwWKegUAywbbUKNi/&#10;bW4fQMSkyeouECq4YIRldX1V6sKGM73iuEuN4BKKhVbgUuoLKaNx6HWchx6Js0MYvE5sh0baQZ+5&#10;3Hcyz7J76XVLvOB0j88OzXF38grievvZm8O2Pjp7+XpZj3fmffOh1M1sWj2BSDilv2P4wWd0qJip&#10;DieyUXQK+JH0q5w9ZhnbWkGeL3KQVSn/01ffAAAA//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAA&#10;AOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAh&#10;ADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAh&#10;AEtsHcAYAgAAJwQAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgA&#10;AAAhAJQbapfbAAAABQEAAA8AAAAAAAAAAAAAAAAAcgQAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAA&#10;BAAEAPMAAAB6BQAAAAA=&#10;"><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p w14:paraId="6BFD21DC" w14:textId="77777777" w:rsidR="00032632" w:rsidRDefault="00032632" w:rsidP="00032632"><w:pPr><w:pStyle w:val="HTMLVorformatiert"/><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>apt-get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> update</w:t></w:r></w:p><w:p w14:paraId="6908EB02" w14:textId="77777777" w:rsidR="00032632" w:rsidRPr="00032632" w:rsidRDefault="00032632" w:rsidP="00032632"><w:pPr><w:pStyle w:val="HTMLVorformatiert"/><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr></w:p><w:p w14:paraId="162C8FC2" w14:textId="0CF419CE" w:rsidR="00032632" w:rsidRPr="00032632" w:rsidRDefault="00032632" w:rsidP="00032632"><w:pPr><w:pStyle w:val="HTMLVorformatiert"/><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>apt-get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>install</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>gz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00032632"><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>garden</w:t></w:r></w:p></w:txbxContent></v:textbox><w10:anchorlock/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$targetParagraph.Range.InsertXML($xml)

# --- 3) Strip the "INTERNAL" classification text box out of every footer --
$sec = $d.Sections.Item(1)
$footers = $sec.Footers
for ($i = 1; $i -le $footers.Count; $i++) {
    $f = $footers.Item($i)
    while ($f.Shapes.Count -gt 0) {
        $f.Shapes.Item(1).Delete()
    }
}

Write-Host "Edits applied."
